$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 33,1
$arr[0,0] = 0.99999999487927749
$arr[1,0] = 0.99637210592182379
$arr[2,0] = 0.98840924688190179
$arr[3,0] = 0.98841340580664194
$arr[4,0] = 0.97517223134352649
$arr[5,0] = 0.94270021312254282
$arr[6,0] = 0.93739399681731106
$arr[7,0] = 0.93002774649147202
$arr[8,0] = 0.92139233372909435
$arr[9,0] = 0.91367012583305529
$arr[10,0] = 0.91254927237045846
$arr[11,0] = 0.91067277277829917
$arr[12,0] = 0.90385783586598456
$arr[13,0] = 0.90173372546861463
$arr[14,0] = 0.90101109764592935
$arr[15,0] = 0.89850452409977644
$arr[16,0] = 0.89479650141142053
$arr[17,0] = 0.89368757537341426
$arr[18,0] = 0.99307867085849844
$arr[19,0] = 0.9663332805429059
$arr[20,0] = 0.95881652540493345
$arr[21,0] = 0.95755200244319738
$arr[22,0] = 0.9813746946202897
$arr[23,0] = 0.96835407482176294
$arr[24,0] = 0.96189711850397663
$arr[25,0] = 0.93904471955350854
$arr[26,0] = 0.93419817431090757
$arr[27,0] = 0.9127253363626926
$arr[28,0] = 0.89745389722883817
$arr[29,0] = 0.8908834530360944
$arr[30,0] = 0.88322955467565833
$arr[31,0] = 0.88155020660701267
$arr[32,0] = 0.881030190481356

$ws.Range("A1:A33").Value = $arr
